# Automatic update of files.
# - Column C ("Förändrad") bumps from 46073 to 46074 for every data row.
# - Rows 4-13 (excluding 9 and 10, whose Beteckning is unique/unchanged) get their
#   identity (Beteckning / Datum / Area) reshuffled among themselves; row 4<->5 swap,
#   and rows 6,7,8,11,12,13 rotate in a 6-cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Bump "Förändrad" (column C) for every row with data (rows 2-13) ----
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# ---- 2. Snapshot the columns that move (A = Beteckning, B = Datum, G = Area) ----
$rows = @(4, 5, 6, 7, 8, 11, 12, 13)
$snapA = @{}
$snapB = @{}
$snapG = @{}
foreach ($r in $rows) {
    $snapA[$r] = $ws.Cells.Item($r, 1).Value()
    $snapB[$r] = $ws.Cells.Item($r, 2).Value()
    $snapG[$r] = $ws.Cells.Item($r, 7).Value()
}

# Mapping: content that WAS at the key row moves to the value row.
$moveTo = @{
    4  = 5
    5  = 4
    6  = 11
    11 = 7
    7  = 12
    12 = 8
    8  = 13
    13 = 6
}

foreach ($src in $moveTo.Keys) {
    $dst = $moveTo[$src]
    $ws.Cells.Item($dst, 1).Value = $snapA[$src]
    $ws.Cells.Item($dst, 2).Value = $snapB[$src]
    $ws.Cells.Item($dst, 7).Value = $snapG[$src]
}

# ---- 3. Rebuild the hyperlink formulas (columns S,T,V,W,X,Y) that embed the
#         Beteckning text; only rows 4 and 5 carry these formulas. ----
function Set-BeteckningLinks($row, $beteckning) {
    $ws.Range("S$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/artfynd/$beteckning artfynd.xlsx"", ""$beteckning"")"
    $ws.Range("T$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/kartor/$beteckning karta.png"", ""$beteckning"")"
    $ws.Range("V$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/klagomål/$beteckning FSC-klagomål.docx"", ""$beteckning"")"
    $ws.Range("W$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/klagomålsmail/$beteckning FSC-klagomål mail.docx"", ""$beteckning"")"
    $ws.Range("X$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/tillsyn/$beteckning tillsynsbegäran.docx"", ""$beteckning"")"
    $ws.Range("Y$row").Formula = "=HYPERLINK(""https://klasma.github.io/Logging_1261/tillsynsmail/$beteckning tillsynsbegäran mail.docx"", ""$beteckning"")"
}

Set-BeteckningLinks 4 $snapA[5]
Set-BeteckningLinks 5 $snapA[4]
